$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.996.04"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "'1.576.43"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'299.33"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").Value = "'0.3743"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").Value = "'0.3549"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").Value = "'49.90"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'1.213"
$ws.Range("E11").Value = "  -4.42%  "
$ws.Range("D12").Value = "'0.07963"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "'21.81"
$ws.Range("E13").Value = "  -5.17%  "
$ws.Range("D14").Value = "'6.415"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "'7.273"
$ws.Range("E15").Value = "  -4.32%  "
$ws.Range("D16").Value = "'0.00001221"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "'1.576.16"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'91.75"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'17.67"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'6.347"
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").Value = "'22.954.28"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "'12.54"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").Value = "'2.372"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'2.810"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").Value = "'20.53"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("D28").Value = "'147.28"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").Value = "'5.160"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'131.56"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'2.331"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("D32").Value = "'6.514"
$ws.Range("E32").Value = "  -6.15%  "
$ws.Range("D33").Value = "'1.752.60"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "'0.9297"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'0.07319"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.08745"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02634"
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("D38").Value = "'9.921"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "'0.2456"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("D40").Value = "'5.972"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("D41").Value = "'1.341"
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("D42").Value = "'0.6831"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").Value = "'11.82"
$ws.Range("E43").Value = "  -7.60%  "
$ws.Range("D44").Value = "'14.65"
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6322"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.958"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.236"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'130.09"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07847"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "'1.185"
$ws.Range("E51").Value = "  +1.45%  "
